$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.272.40"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "4.033.91"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'530.60"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'151.68"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Value = "'0.703"
$ws.Range("E7").Value = "  +12.88%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.752"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "'51.14"
$ws.Range("E11").Value = "  +8.03%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "'0.0000327"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").Value = "'10.79"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "4.679.33"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "4.033.75"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'14.12"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "72.213.20"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'434.10"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "'98.15"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").Value = "'3.51"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'4.22"
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("D25").Value = "'14.39"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'11.21"
$ws.Range("E26").Value = "  -7.50%  "
$ws.Range("D27").Value = "'10.79"
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'5.86"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'3.66"
$ws.Range("E29").Value = "  +18.41%  "
$ws.Range("D30").Value = "'36.84"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "'7.49"
$ws.Range("E31").Value = "  +7.63%  "
$ws.Range("D32").Value = "'13.51"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'0.132"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").Value = "'681.94"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "'48.36"
$ws.Range("E35").Value = "  +18.80%  "
$ws.Range("D36").Value = "'65.71"
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.152"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0828"
$ws.Range("E39").Value = "  -8.63%  "
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -7.64%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'3.36"
$ws.Range("E41").Value = "  +7.95%  "
$ws.Range("B42").Value = "Dai"
$ws.Range("C42").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'10.36"
$ws.Range("E45").Value = "  +12.98%  "
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").Value = "'2.68"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("D48").Value = "'3.40"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("D49").Value = "'3.03"
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'3.28"
$ws.Range("E51").Value = "  -1.61%  "
